# carpentry_mrp_import / report_mrp_component_unknown.xlsx
# - Insert a new "section title" row above the column-header row.
# - Re-word the report title in A1.
# - Style + size the new section-title row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row 6 (pushes the old header row 6 -> row 7).
$ws.Rows(6).Insert()

# New "Section title" cell in the freshly inserted row.
$ws.Range("A6").Value = "Section title"
$ws.Range("A6").Font.Size = 14
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Rows(6).RowHeight = 18

# Reword the main report title.
$ws.Range("A1").Value = "Report for the import of manufacturing Components"

# Leave the freshly-inserted rows selected, matching the authored state.
[void]$ws.Range("A6:XFD7").Select()
